$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.296.54'
$cell.Style = $origStyle
$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.75%  '
$cell.Style = $origStyle
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.664.07'
$cell.Style = $origStyle
$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.53%  '
$cell.Style = $origStyle
$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.84%  '
$cell.Style = $origStyle
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '219.07'
$cell.Style = $origStyle
$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.55%  '
$cell.Style = $origStyle
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5341'
$cell.Style = $origStyle
$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.37%  '
$cell.Style = $origStyle
$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.80%  '
$cell.Style = $origStyle
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.2651'
$cell.Style = $origStyle
$cell = $ws.Range("E8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.51%  '
$cell.Style = $origStyle
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06413'
$cell.Style = $origStyle
$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.02%  '
$cell.Style = $origStyle
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '20.61'
$cell.Style = $origStyle
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.07825'
$cell.Style = $origStyle
$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.40%  '
$cell.Style = $origStyle
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.575'
$cell.Style = $origStyle
$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.46%  '
$cell.Style = $origStyle
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.673.72'
$cell.Style = $origStyle
$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.98%  '
$cell.Style = $origStyle
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.892.31'
$cell.Style = $origStyle
$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.50%  '
$cell.Style = $origStyle
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5527'
$cell.Style = $origStyle
$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.62%  '
$cell.Style = $origStyle
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0₅8215'
$cell.Style = $origStyle
$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.09%  '
$cell.Style = $origStyle
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '65.68'
$cell.Style = $origStyle
$cell = $ws.Range("E17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.37%  '
$cell.Style = $origStyle
$cell = $ws.Range("B18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'WrappedBTC'
$cell.Style = $origStyle
$cell = $ws.Range("C18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$cell.Style = $origStyle
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.304.25'
$cell.Style = $origStyle
$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.69%  '
$cell.Style = $origStyle
$cell = $ws.Range("B19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Dai'
$cell.Style = $origStyle
$cell = $ws.Range("C19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell.Style = $origStyle
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.010'
$cell.Style = $origStyle
$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.84%  '
$cell.Style = $origStyle
$cell = $ws.Range("B20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Uniswap'
$cell.Style = $origStyle
$cell = $ws.Range("C20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell.Style = $origStyle
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.692'
$cell.Style = $origStyle
$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.33%  '
$cell.Style = $origStyle
$cell = $ws.Range("B21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'BitcoinCash'
$cell.Style = $origStyle
$cell = $ws.Range("C21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell.Style = $origStyle
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '193.48'
$cell.Style = $origStyle
$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.35%  '
$cell.Style = $origStyle
$cell = $ws.Range("B22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Avalanche'
$cell.Style = $origStyle
$cell = $ws.Range("C22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell.Style = $origStyle
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.21'
$cell.Style = $origStyle
$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.40%  '
$cell.Style = $origStyle
$cell = $ws.Range("B23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Chainlink'
$cell.Style = $origStyle
$cell = $ws.Range("C23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell.Style = $origStyle
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.040'
$cell.Style = $origStyle
$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.12%  '
$cell.Style = $origStyle
$cell = $ws.Range("B24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'BinanceUSD'
$cell.Style = $origStyle
$cell = $ws.Range("C24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell.Style = $origStyle
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.011'
$cell.Style = $origStyle
$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.81%  '
$cell.Style = $origStyle
$cell = $ws.Range("B25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Monero'
$cell.Style = $origStyle
$cell = $ws.Range("C25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell.Style = $origStyle
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '146.00'
$cell.Style = $origStyle
$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.28%  '
$cell.Style = $origStyle
$cell = $ws.Range("B26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Stellar'
$cell.Style = $origStyle
$cell = $ws.Range("C26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell.Style = $origStyle
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.1235'
$cell.Style = $origStyle
$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.13%  '
$cell.Style = $origStyle
$cell = $ws.Range("B27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Cosmos'
$cell.Style = $origStyle
$cell = $ws.Range("C27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell.Style = $origStyle
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.197'
$cell.Style = $origStyle
$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.76%  '
$cell.Style = $origStyle
$cell = $ws.Range("B28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'EthereumClassic'
$cell.Style = $origStyle
$cell = $ws.Range("C28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell.Style = $origStyle
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '16.12'
$cell.Style = $origStyle
$cell = $ws.Range("E28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.18%  '
$cell.Style = $origStyle
$cell = $ws.Range("B29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Toncoin'
$cell.Style = $origStyle
$cell = $ws.Range("C29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell.Style = $origStyle
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.483'
$cell.Style = $origStyle
$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.80%  '
$cell.Style = $origStyle
$cell = $ws.Range("B30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Hedera'
$cell.Style = $origStyle
$cell = $ws.Range("C30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell.Style = $origStyle
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.05850'
$cell.Style = $origStyle
$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.72%  '
$cell.Style = $origStyle
$cell = $ws.Range("B31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'PancakeSwap'
$cell.Style = $origStyle
$cell = $ws.Range("C31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell.Style = $origStyle
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.282'
$cell.Style = $origStyle
$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.54%  '
$cell.Style = $origStyle
$cell = $ws.Range("B32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'InternetComputer(DFINITY)'
$cell.Style = $origStyle
$cell = $ws.Range("C32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell.Style = $origStyle
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.623'
$cell.Style = $origStyle
$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.94%  '
$cell.Style = $origStyle
$cell = $ws.Range("B33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Filecoin'
$cell.Style = $origStyle
$cell = $ws.Range("C33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell.Style = $origStyle
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.283'
$cell.Style = $origStyle
$cell = $ws.Range("E33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.55%  '
$cell.Style = $origStyle
$cell = $ws.Range("B34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'LidoDAOToken'
$cell.Style = $origStyle
$cell = $ws.Range("C34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell.Style = $origStyle
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.612'
$cell.Style = $origStyle
$cell = $ws.Range("B35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'ARBITRUM'
$cell.Style = $origStyle
$cell = $ws.Range("C35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell.Style = $origStyle
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.9631'
$cell.Style = $origStyle
$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.09%  '
$cell.Style = $origStyle
$cell = $ws.Range("B36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'MXToken'
$cell.Style = $origStyle
$cell = $ws.Range("C36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell.Style = $origStyle
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.826'
$cell.Style = $origStyle
$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.49%  '
$cell.Style = $origStyle
$cell = $ws.Range("B37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'HuobiToken'
$cell.Style = $origStyle
$cell = $ws.Range("C37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell.Style = $origStyle
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.417'
$cell.Style = $origStyle
$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.22%  '
$cell.Style = $origStyle
$cell = $ws.Range("B38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'ImmutableX'
$cell.Style = $origStyle
$cell = $ws.Range("C38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell.Style = $origStyle
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5803'
$cell.Style = $origStyle
$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.78%  '
$cell.Style = $origStyle
$cell = $ws.Range("B39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'VeChain'
$cell.Style = $origStyle
$cell = $ws.Range("C39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell.Style = $origStyle
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.01609'
$cell.Style = $origStyle
$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.49%  '
$cell.Style = $origStyle
$cell = $ws.Range("B40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'TrustWalletToken'
$cell.Style = $origStyle
$cell = $ws.Range("C40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell.Style = $origStyle
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.8673'
$cell.Style = $origStyle
$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.18%  '
$cell.Style = $origStyle
$cell = $ws.Range("B41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'FraxShare'
$cell.Style = $origStyle
$cell = $ws.Range("C41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell.Style = $origStyle
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.890'
$cell.Style = $origStyle
$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.38%  '
$cell.Style = $origStyle
$cell = $ws.Range("B42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Maker'
$cell.Style = $origStyle
$cell = $ws.Range("C42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$cell.Style = $origStyle
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.051.19'
$cell.Style = $origStyle
$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.10%  '
$cell.Style = $origStyle
$cell = $ws.Range("B43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'PaxDollar'
$cell.Style = $origStyle
$cell = $ws.Range("C43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell.Style = $origStyle
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.010'
$cell.Style = $origStyle
$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.76%  '
$cell.Style = $origStyle
$cell = $ws.Range("B44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Quant'
$cell.Style = $origStyle
$cell = $ws.Range("C44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell.Style = $origStyle
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '104.52'
$cell.Style = $origStyle
$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.86%  '
$cell.Style = $origStyle
$cell = $ws.Range("B45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'RocketPoolETH'
$cell.Style = $origStyle
$cell = $ws.Range("C45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$cell.Style = $origStyle
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.803.06'
$cell.Style = $origStyle
$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.27%  '
$cell.Style = $origStyle
$cell = $ws.Range("B46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Aave'
$cell.Style = $origStyle
$cell = $ws.Range("C46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell.Style = $origStyle
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '57.72'
$cell.Style = $origStyle
$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.02%  '
$cell.Style = $origStyle
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0₈107'
$cell.Style = $origStyle
$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.57%  '
$cell.Style = $origStyle
$cell = $ws.Range("B48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Frax'
$cell.Style = $origStyle
$cell = $ws.Range("C48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$cell.Style = $origStyle
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.012'
$cell.Style = $origStyle
$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.01%  '
$cell.Style = $origStyle
$cell = $ws.Range("B49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Mantle'
$cell.Style = $origStyle
$cell = $ws.Range("C49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$cell.Style = $origStyle
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.4383'
$cell.Style = $origStyle
$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.96%  '
$cell.Style = $origStyle
$cell = $ws.Range("B50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'EnergySwap'
$cell.Style = $origStyle
$cell = $ws.Range("C50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell.Style = $origStyle
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.005'
$cell.Style = $origStyle
$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.31%  '
$cell.Style = $origStyle
$cell = $ws.Range("B51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Cronos'
$cell.Style = $origStyle
$cell = $ws.Range("C51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell.Style = $origStyle
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.05166'
$cell.Style = $origStyle
$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.28%  '
$cell.Style = $origStyle
